$d = $word.ActiveDocument

# The first paragraph in the document holds the hidden bookmark-style
# "**ID__...__ID**" placeholder text.
$para1 = $d.Paragraphs.Item(1)

# 1) Give the paragraph a border (top/left/bottom/right, 5-twip space,
#    no visible line) -- matches <w:pBdr><w:top w:space="5"/>...</w:pBdr>
$pBdr = $para1.Range.ParagraphFormat.Borders
$pBdr.DistanceFromTop = 5
$pBdr.DistanceFromBottom = 5
$pBdr.DistanceFromLeft = 5
$pBdr.DistanceFromRight = 5

# 2) Bump the left indent from 120 to 225 twips (COM works in points,
#    1 pt = 20 twips -> 225/20 = 11.25 pt).
$para1.Range.ParagraphFormat.LeftIndent = 11.25

# 3) Update the placeholder id text and drop the trailing space run that
#    used to follow it (the old text + trailing space is replaced by the
#    new id text with no trailing space, removing the second run).
$para1.Range.Find.Execute("**ID__AFFARS_pgi_5301_topic_23__ID** ", $true, $false, $false, $false, $false,
                           $true, 1, $false, "**ID__AFFARS_AFMC_PGI_5301_290__ID**", 2)
